# "Generate Report for Handback"
#
# The handback transform failed for the 92d54326-... item (row 7 of each
# data sheet). Update the status everywhere it is surfaced, and record the
# handback/handoff filename-mismatch error detail for both locales.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet: File Name row for 92d54326..., zh-cn + de-de status cells (E7, F7)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E7").Value = $newStatus
$wsOverview.Range("F7").Value = $newStatus

# --- zh-cn sheet: Status (C7) + Error Detail (P7) for the same row
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C7").Value = $newStatus
$wsZh.Range("P7").Value = "Handback file name: osbz2wh4.t5m is different with handoff file name: 92d54326-6331-4bda-b9a5-3e56ddfec95b.22c6572d19ca3ab0bcb162e8ac9727da551aadd7.zh-cn."
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: Status (C7) + Error Detail (P7) for the same row
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C7").Value = $newStatus
$wsDe.Range("P7").Value = "Handback file name: osbz2wh4.t5m is different with handoff file name: 92d54326-6331-4bda-b9a5-3e56ddfec95b.22c6572d19ca3ab0bcb162e8ac9727da551aadd7.de-de."
$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664
